# Change request to remove some sites #99
#
# - Removes Tsakane, Katlehong and the Ekurhuleni/Bojanala "Mobile" sites
#   from the `sites` choice list.
# - Adds a new "Bapong Clinic" (bapong) site under the Bojanala district
#   (replacing the old "mobile-b" row).
# - Extends the AE code regex constraint to also accept lower-case codes.
# - Updates sheet selection / active tab to match the author's last
#   on-screen state.

$wb = $excel.ActiveWorkbook
$surveySheet  = $wb.Worksheets.Item("survey")
$choicesSheet = $wb.Worksheets.Item("choices")

# --- choices: drop Tsakane / Katlehong / Mobile (ekurhuleni) rows -----------
# Original rows (1-based):
#   23 winnie_mandela / Winnie Mandela Clinic / ekurhuleni   (kept)
#   24 tsakane        / Tsakane Clinic        / ekurhuleni   (removed)
#   25 katlehong      / Katlehong Clinic      / ekurhuleni   (removed)
#   26 mobile-e       / Mobile                / ekurhuleni   (removed)
#   27 bafokeng       / Bafokeng              / bojanala     (kept, shifts up)
#   28 letlhabile     / Letlhabile            / bojanala     (kept, shifts up)
#   29 mogwase        / Mogwase Clinic        / bojanala     (kept, shifts up)
#   30 mobile-b       / Mobile                / bojanala     (replaced by bapong)
$choicesSheet.Range("A24:D26").EntireRow.Delete()

# The old "mobile-b" site (now shifted up to row 27) becomes the new
# "Bapong Clinic" site, still under the Bojanala district.
$choicesSheet.Range("B27").Value = "bapong"
$choicesSheet.Range("C27").Value = "Bapong Clinic"

# --- survey: widen the ae_code regex constraint to also allow lower case ----
$surveySheet.Range("G28").Value = 'regex(., ''^(AN|BL|IN|OT|PA|SD|SX|WD|OA)-[A-C]$'') or regex(., ''^(an|bl|in|ot|pa|sd|sx|wd|oa)-[a-c]$'')'

# --- view state: "choices" tab becomes the active/selected sheet -----------
$surveySheet.Range("G28").Select()
$choicesSheet.Activate()
$choicesSheet.Range("D33").Select()
